$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect it so the locked header cell can be edited.
$ws.Unprotect()
Start-Sleep -Milliseconds 100

# Shorten the header text (used by cell A1 and by Table2's first column header).
$ws.Range("A1").Value = "Lecture start time"

# Bump the header cell's font size from 6pt to 10pt.
$ws.Range("A1").Font.Size = 10

# Restore sheet protection.
$ws.Protect()
Start-Sleep -Milliseconds 100

$wb.Save()
Write-Host "Saved workbook"
